$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-09-15")
$src.Copy([System.Reflection.Missing]::Value, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "2025-09-22"

$ws.Range("B2").Value = "ワンパンマン"
$ws.Range("C2").Value = "原作/ＯＮＥ 作画/村田雄介"
$ws.Range("D2").Value = "211撃目"
$ws.Range("B3").Value = "異種族レビュアーズ"
$ws.Range("C3").Value = "天原(原作) masha(作画)"
$ws.Range("D3").Value = "第87話"
$ws.Range("B4").Value = "新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる"
$ws.Range("C4").Value = "漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea"
$ws.Range("D4").Value = "第71話"
$ws.Range("B5").Value = "異世界おじさん"
$ws.Range("C5").Value = "殆ど死んでいる(著者)"
$ws.Range("D5").Value = "コミックス第14巻発売告知"
$ws.Range("B6").Value = "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～"
$ws.Range("C6").Value = "zunta(作画) はらわたさいぞう(原作)"
$ws.Range("D6").Value = "第32話：思考を奪う③"
$ws.Range("B7").Value = "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―"
$ws.Range("C7").Value = "光永康則"
$ws.Range("D7").Value = "第６９話『岩鬼停止』③"
$ws.Range("B8").Value = "勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが"
$ws.Range("C8").Value = "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)"
$ws.Range("D8").Value = "第5話 前編"
$ws.Range("B9").Value = "クセ強彼女は床にいざなう"
$ws.Range("C9").Value = "須河篤志(著者)"
$ws.Range("D9").Value = "第15話後半"
$ws.Range("B10").Value = "追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。"
$ws.Range("C10").Value = "六志麻あさ 業務用餅 kisui"
$ws.Range("D10").Value = "第７２話"
$ws.Range("B11").Value = "美人女上司滝沢さん"
$ws.Range("C11").Value = "やんBARU(著者)"
$ws.Range("D11").Value = "第203話"
$ws.Range("B12").Value = "実は俺、最強でした？"
$ws.Range("C12").Value = "原作：澄守 彩 漫画：高橋 愛"
$ws.Range("D12").Value = "第126話　王家の立場"
$ws.Range("B13").Value = "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。"
$ws.Range("C13").Value = "マツモトケンゴ"
$ws.Range("D13").Value = "第64話 更衣室の戦いが始まった（２）"
$ws.Range("B14").Value = "帰ってください！ 阿久津さん"
$ws.Range("C14").Value = "長岡太一(著者)"
$ws.Range("D14").Value = "第196話"
$ws.Range("B15").Value = "元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～"
$ws.Range("C15").Value = "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)"
$ws.Range("D15").Value = "第78話その2"
$ws.Range("B16").Value = "異世界魔王と召喚少女の奴隷魔術"
$ws.Range("C16").Value = "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大"
$ws.Range("D16").Value = "第128話　レムと話してみる（後編）"
$ws.Range("B17").Value = "ダークサモナーとデキている"
$ws.Range("C17").Value = "車王(著者)"
$ws.Range("D17").Value = "第76話"
$ws.Range("B18").Value = "蜘蛛ですが、なにか？"
$ws.Range("C18").Value = "かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)"
$ws.Range("D18").Value = "第76話その1"
$ws.Range("B19").Value = "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～"
$ws.Range("C19").Value = "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり"
$ws.Range("D19").Value = "第５２話　暴走を止める器用貧乏（４）"
$ws.Range("B20").Value = "まんきつしたい常連さん"
$ws.Range("C20").Value = "しんみりん(著者)"
$ws.Range("D20").Value = "第48話前編"
$ws.Range("B21").Value = "淫獄団地"
$ws.Range("C21").Value = "搾精研究所(原作) 丈山雄為(漫画)"
$ws.Range("D21").Value = "第50話（後編）"
$ws.Range("B22").Value = "Ｓ級ギルドを追放されたけど、実は俺だけドラゴンの言葉がわかるので、気付いたときには竜騎士の頂点を極めてました。"
$ws.Range("C22").Value = "ひそな(漫画) 三木なずな(原作) 白狼(キャラクター原案)"
$ws.Range("D22").Value = "第39話-2"
$ws.Range("B23").Value = "異世界のんびり農家"
$ws.Range("C23").Value = "剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)"
$ws.Range("D23").Value = "第305話"
$ws.Range("B24").Value = "バキ外伝 烈海王は異世界転生しても一向にかまわんッッ"
$ws.Range("C24").Value = "板垣恵介 猪原賽 陸井栄史"
$ws.Range("D24").Value = "第80話　先に行く"
$ws.Range("B25").Value = "リビルドワールド"
$ws.Range("C25").Value = "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)"
$ws.Range("D25").Value = "第73話②"
$ws.Range("B26").Value = "戸崎さんは僕にだけ冷たい"
$ws.Range("C26").Value = "saku(著者)"
$ws.Range("D26").Value = "第29話-1"
$ws.Range("B27").Value = "濁る瞳で何を願う ハイセルク戦記"
$ws.Range("C27").Value = "トルトネン 創-taro 斎藤八呑"
$ws.Range("D27").Value = "第34話 ダンデューグ城へようこそ"
$ws.Range("B28").Value = "ニチアサ好きのオタクが悪役生徒に転生した結果、破滅フラグが崩壊していく件について"
$ws.Range("C28").Value = "烏丸英（原作） どんぐりす（漫画）"
$ws.Range("D28").Value = "第14話（後編）急襲…事件の始まり"
$ws.Range("B29").Value = "モブ高生の俺でも冒険者になればリア充になれますか？"
$ws.Range("C29").Value = "原作：百均 漫画：さぎやまれん キャラクター原案：hai"
$ws.Range("D29").Value = "第31話"
$ws.Range("B30").Value = "君のラブを見せてくれ！"
$ws.Range("C30").Value = "リムコロ(著者)"
$ws.Range("D30").Value = "コミックス第⑤巻発売告知"
$ws.Range("B31").Value = "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～"
$ws.Range("C31").Value = "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき"
$ws.Range("D31").Value = "第34話 独身貴族は礼の品を贈る（1）"
$ws.Range("B32").Value = "ルパン三世 異世界の姫君（ネイバーワールドプリンセス）"
$ws.Range("C32").Value = "モンキー・パンチ／エム・ピー・ワークス 内々けやき 佐伯庸介 白狼"
$ws.Range("D32").Value = "第110話：王女に贈る子守歌"
$ws.Range("B33").Value = "「おかえり、パパ」"
$ws.Range("C33").Value = "蝉丸"
$ws.Range("D33").Value = "第28話　帰宅"
$ws.Range("B34").Value = "よくわからないけれど異世界に転生していたようです"
$ws.Range("C34").Value = "内々けやき あし カオミン"
$ws.Range("D34").Value = "第139話 よくわからないけれどズゥゥゥンときたみたいです（２）"
$ws.Range("B35").Value = "聖者無双"
$ws.Range("C35").Value = "漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime"
$ws.Range("D35").Value = "第92話　龍と精霊の信仰（後半）"
$ws.Range("B36").Value = "落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～"
$ws.Range("C36").Value = "村上よしゆき 茨木野 あるてら"
$ws.Range("D36").Value = "第４２話　勇者、六邪神将相手に舐めプしてたら、ピンチになる（３）"
$ws.Range("B37").Value = "ライドンキング"
$ws.Range("C37").Value = "馬場康誌"
$ws.Range("D37").Value = "第83話 大統領と龍の闇卵（前編）"
$ws.Range("B38").Value = "バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～"
$ws.Range("C38").Value = "板垣恵介 林たかあき"
$ws.Range("D38").Value = "第54話 イバラの帰り道"
$ws.Range("B39").Value = "願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜"
$ws.Range("C39").Value = "ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)"
$ws.Range("D39").Value = "第6話-2：火蓮の剣"
$ws.Range("B40").Value = "小林さんちのメイドラゴン"
$ws.Range("C40").Value = "クール教信者"
$ws.Range("D40").Value = "第150話"
$ws.Range("B41").Value = "治癒魔法の間違った使い方 ~戦場を駆ける回復要員~"
$ws.Range("C41").Value = "九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)"
$ws.Range("D41").Value = "第82話その4"
$ws.Range("B42").Value = "宇崎ちゃんは遊びたい！"
$ws.Range("C42").Value = "丈(著者)"
$ws.Range("D42").Value = "第127話"
$ws.Range("B43").Value = "ダメ人間の愛しかた"
$ws.Range("C43").Value = "岩葉(著者)"
$ws.Range("D43").Value = "第20話前編　ダメ人間と3人暮らしの彼女"
$ws.Range("B44").Value = "婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版"
$ws.Range("C44").Value = "漫画/すたひろ 原作/Y.A"
$ws.Range("D44").Value = "chapter70【37話①】"
$ws.Range("B45").Value = "理想のヒモ生活"
$ws.Range("C45").Value = "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)"
$ws.Range("D45").Value = "第87話　その3"
$ws.Range("B46").Value = "彼女にしたい女子一位、の隣で見つけたあまりちゃん"
$ws.Range("C46").Value = "寝巻ネルゾ(漫画) 裕時悠示(原作) たん旦(キャラクター原案)"
$ws.Range("D46").Value = "第5話②「4人で遊園地！」"
$ws.Range("B47").Value = "10年ぶりに再会したクソガキは清純美少女JKに成長していた"
$ws.Range("C47").Value = "緑青黒羽（漫画） 館西夕木（原作） ひげ猫（キャラクター原案）"
$ws.Range("D47").Value = "第6話　美少女二人（後編）"
$ws.Range("B48").Value = "生徒会役員共"
$ws.Range("C48").Value = "氏家ト全"
$ws.Range("D48").Value = "#414"
$ws.Range("B49").Value = "俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。"
$ws.Range("C49").Value = "狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)"
$ws.Range("D49").Value = "第25話-2"
$ws.Range("B50").Value = "ふかふかダンジョン攻略記～俺の異世界転生冒険譚～"
$ws.Range("C50").Value = "KAKERU"
$ws.Range("D50").Value = "第68話「東アイギス2」（後半）"
$ws.Range("B51").Value = "賢者の孫"
$ws.Range("C51").Value = "緒方俊輔(漫画) 吉岡剛(原作) 菊池政治(キャラクター原案)"
$ws.Range("D51").Value = "第95話-3"
